# Add a new "status" column (column CK) to the feature table on Sheet1.
# This mirrors the existing "page_rank" column (CJ): same header style,
# and a value of 0 for every data row (rows 2-175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastHeaderCol = "CJ"
$newCol = "CK"

$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 2) {
    $lastRow = 175
}

# Copy the formatting (bold font + border + alignment) of the last header
# cell onto the new header cell so the new column looks like the others.
$ws.Range($lastHeaderCol + "1").Copy()
$ws.Range($newCol + "1").PasteSpecial(-4122)

# Set the header text for the new column.
$ws.Range($newCol + "1").Value = "status"

# Fill the new column with 0 for every existing data row.
$dataRange = $ws.Range($newCol + "2:" + $newCol + $lastRow)
$dataRange.Value = 0
